$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8378782272338867
$ws.Range("B1").Value = 1.254250645637512
$ws.Range("C1").Value = 2.292247772216797
$ws.Range("D1").Value = 2.413295745849609
$ws.Range("E1").Value = 1.966107845306396
